$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / unambiguous cell updates (safe to set directly) ---
$ws.Range("D2").Value = "29.405.71"
$ws.Range("D3").Value = "1.878.38"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +1.50%  "
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("E11").Value = "  -3.15%  "
$ws.Range("D12").Value = "1.887.54"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("E13").Value = "  +3.88%  "
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("E16").Value = "  +3.37%  "
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "29.418.25"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("E19").Value = "  +5.14%  "
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("D21").Value = "2.135.69"
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("E23").Value = "  -1.33%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("E28").Value = "  +2.43%  "
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("E32").Value = "  +1.22%  "
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("E35").Value = "  +1.24%  "
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("D39").Value = "1.278.45"
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("E41").Value = "  -2.49%  "
$ws.Range("E42").Value = "  +1.48%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("E43").Value = "  +0.60%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("E44").Value = "  +1.62%  "
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("D47").Value = "2.030.65"
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("E51").Value = "  -0.42%  "

# --- Numeric-looking text cells: force text via NumberFormat, then restore default style ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7171"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.71"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07975"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3148"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.94"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "94.76"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.235"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7079"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.417"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008430"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "253.05"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.35"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.681"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1581"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.077"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.33"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.98"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.422"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.319"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.223"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7582"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.703"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01885"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.760"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.407"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9086"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "111.90"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "74.32"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5208"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.530"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4341"

# Restore the default (unstyled) format on the numeric-looking text cells so they
# match the original style (no explicit style index), mirroring a scratch blank cell.
$scratch = $ws.Cells.Item(200, 200)
$scratch.Copy()
$restoreList = @("D5","D6","D8","D9","D10","D13","D14","D15","D16","D17","D19","D20","D22","D23","D24","D25","D26","D27","D28","D30","D31","D32","D35","D37","D38","D40","D41","D42","D43","D44","D49","D50","D51")
foreach ($ref in $restoreList) {
    $ws.Range($ref).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false
